# Jobs_Positions_de.xlsx update
# - Translate area_title values ("Manual"/"Small Parts"/"High Rack") to German
#   ("Manuell"/"Kleinteile"/"Hochregal") throughout column C of the data table,
#   except for the "Small Parts" occurrences in rows 41-43 and 57-59 which were
#   left untranslated in the source edit.
# - Apply an AutoFilter over the data range A1:F61 (adds the hidden
#   _xlnm._FilterDatabase defined name as a side effect).
# - Update the active selection to C36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Translate "Manual" -> "Manuell" ---
$ws.Range("C2:C4").Value = "Manuell"
$ws.Range("C17").Value = "Manuell"
$ws.Range("C25:C31").Value = "Manuell"
$ws.Range("C40").Value = "Manuell"
$ws.Range("C52:C56").Value = "Manuell"

# --- Translate "Small Parts" -> "Kleinteile" ---
# (rows 41-43 and 57-59 intentionally left as "Small Parts")
$ws.Range("C5:C8").Value = "Kleinteile"
$ws.Range("C18:C20").Value = "Kleinteile"
$ws.Range("C32:C35").Value = "Kleinteile"

# --- Translate "High Rack" -> "Hochregal" ---
$ws.Range("C9:C16").Value = "Hochregal"
$ws.Range("C21:C24").Value = "Hochregal"
$ws.Range("C36:C39").Value = "Hochregal"
$ws.Range("C44:C51").Value = "Hochregal"
$ws.Range("C60:C61").Value = "Hochregal"

# --- Apply AutoFilter over the data range ---
$ws.Range("A1:F61").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=data!`$A`$1:`$F`$61")
$filterName.Visible = $false

# --- Update active selection ---
$ws.Range("C36").Select()

Write-Host "done"
